$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings need to be
# forced to Text format first, otherwise Excel auto-converts them to
# numbers (losing the original text formatting / trailing zeros etc.),
# then the style is reset back to Normal so no stray formatting remains.

$ws.Range('D2').Value = '25.794.34'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.636.64'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.86'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.65%  '
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.25'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.637.80'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').Value = '1.863.18'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.555'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.46%  '
$ws.Range('D16').Value = '0.0₃0777'
$ws.Range('E16').Value = '  +2.29%  '
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = '25.825.77'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('E20').Value = '  +2.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '194.45'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.15'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.96%  '
$ws.Range('E25').Value = '  -2.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.13'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('E27').Value = '  -5.21%  '
$ws.Range('E28').Value = '  +0.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.58'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.80%  '
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0493'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.41%  '
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.25'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.10%  '
$ws.Range('E34').Value = '  +2.02%  '
$ws.Range('E35').Value = '  +0.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.898'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.551'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('D39').Value = '1.113.13'
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.57'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '98.99'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.18%  '
$ws.Range('E44').Value = '  -0.23%  '
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.58'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.51'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +12.31%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.74'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.418'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.36%  '
$ws.Range('E50').Value = '  -0.42%  '
$ws.Range('E51').Value = '  -0.44%  '
